$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the "target_market" text column (C2:C266, currently all the
# shared-string "18-25") with a computed formula =74500000*0.17.
$ws.Range("C2:C266").Formula = "=74500000*0.17"

# Move the active selection to C1 (was C20).
$ws.Range("C1").Select()
